$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1410.8541
$ws.Range("I15").Value = 1410.8541
$ws.Range("K15").Value = 4232.5623
$ws.Range("M15").Value = -4063.5623

$ws.Range("H45").Value = 1017
$ws.Range("I45").Value = 1017
$ws.Range("K45").Value = 3051
$ws.Range("M45").Value = -2859

$ws.Range("H74").Value = 2504.2173
$ws.Range("I74").Value = 2527.1365
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 2527.1365
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -1591.1365
$ws.Range("N74").Value = -3872

$ws.Range("H77").Value = 2504.2173
$ws.Range("I77").Value = 2527.1365
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 12635.6825
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -7955.682500000001
$ws.Range("N77").Value = -19360

$ws.Range("H137").Value = 1320.711
$ws.Range("I137").Value = 1205.9688
$ws.Range("J137").Value = 1603.1538
$ws.Range("K137").Value = 3617.9064
$ws.Range("L137").Value = 4809.4614
$ws.Range("M137").Value = -1067.9064
$ws.Range("N137").Value = -9909.4614

$ws.Range("H141").Value = 2796.08
$ws.Range("I141").Value = 1247.2927
$ws.Range("J141").Value = 9851.666999999999
$ws.Range("K141").Value = 3741.8781
$ws.Range("L141").Value = 29555.001
$ws.Range("M141").Value = 1438.1219
$ws.Range("N141").Value = -39915.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1495.9231
$ws.Range("I31").Value = 1070.5
$ws.Range("J31").Value = 1721.1471
$ws.Range("K31").Value = 1070.5
$ws.Range("L31").Value = 1721.1471
$ws.Range("M31").Value = -775.5
$ws.Range("N31").Value = -2311.1471

$ws.Range("H34").Value = 1495.9231
$ws.Range("I34").Value = 1070.5
$ws.Range("J34").Value = 1721.1471
$ws.Range("K34").Value = 1070.5
$ws.Range("L34").Value = 1721.1471
$ws.Range("M34").Value = -868.5
$ws.Range("N34").Value = -2125.1471

$ws.Range("H58").Value = 13890045
$ws.Range("I58").Value = 19609102
$ws.Range("J58").Value = 906
$ws.Range("K58").Value = 19609102
$ws.Range("L58").Value = 906
$ws.Range("M58").Value = -19608899
$ws.Range("N58").Value = -1312

$ws.Range("H80").Value = 23800
$ws.Range("J80").Value = 23800
$ws.Range("L80").Value = 23800
$ws.Range("N80").Value = -26046

$ws.Range("H82").Value = 20390.334
$ws.Range("J82").Value = 20390.334
$ws.Range("L82").Value = 20390.334
$ws.Range("N82").Value = -21112.334

$ws.Range("H83").Value = 23800
$ws.Range("J83").Value = 23800
$ws.Range("L83").Value = 71400
$ws.Range("N83").Value = -82632

$ws.Range("H85").Value = 20390.334
$ws.Range("J85").Value = 20390.334
$ws.Range("L85").Value = 20390.334
$ws.Range("N85").Value = -22886.334

$ws.Range("H132").Value = 7937564
$ws.Range("I132").Value = 1250.2273
$ws.Range("J132").Value = 16667509
$ws.Range("K132").Value = 3750.6819
$ws.Range("L132").Value = 50002527
$ws.Range("M132").Value = -1220.6819
$ws.Range("N132").Value = -50007587

$ws.Range("H134").Value = 1190.0454
$ws.Range("I134").Value = 952.6
$ws.Range("J134").Value = 1698.8572
$ws.Range("K134").Value = 2857.8
$ws.Range("L134").Value = 5096.571599999999
$ws.Range("M134").Value = -322.8000000000002
$ws.Range("N134").Value = -10166.5716

$ws.Range("H136").Value = 13890045
$ws.Range("I136").Value = 19609102
$ws.Range("J136").Value = 906
$ws.Range("K136").Value = 58827306
$ws.Range("L136").Value = 2718
$ws.Range("M136").Value = -58824756
$ws.Range("N136").Value = -7818

$ws.Range("H141").Value = 28980.125
$ws.Range("J141").Value = 28980.125
$ws.Range("L141").Value = 28980.125
$ws.Range("N141").Value = -39340.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20000044
$ws.Range("I2").Value = 48.333332
$ws.Range("J2").Value = 200000000
$ws.Range("K2").Value = 289.999992
$ws.Range("L2").Value = 1200000000
$ws.Range("M2").Value = -176.999992
$ws.Range("N2").Value = -1200000226

$ws.Range("H74").Value = 500000500
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -1939

$ws.Range("H77").Value = 500000500
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 9000
$ws.Range("M77").Value = -3696

$ws.Range("H131").Value = 928.1
$ws.Range("J131").Value = 952.1489
$ws.Range("L131").Value = 2856.4467
$ws.Range("N131").Value = -12936.4467

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5266217
$ws.Range("J80").Value = 33334900
$ws.Range("L80").Value = 33334900
$ws.Range("N80").Value = -33336896

$ws.Range("H83").Value = 5266217
$ws.Range("J83").Value = 33334900
$ws.Range("L83").Value = 166674500
$ws.Range("N83").Value = -166684484

$ws.Range("H97").Value = 981.4643
$ws.Range("I97").Value = 721.5217
$ws.Range("J97").Value = 2177.2
$ws.Range("K97").Value = 721.5217
$ws.Range("L97").Value = 2177.2
$ws.Range("M97").Value = -225.5217
$ws.Range("N97").Value = -3169.2

$ws.Range("H105").Value = 49800
$ws.Range("J105").Value = 49800
$ws.Range("L105").Value = 49800
$ws.Range("N105").Value = -56788

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1850
$ws.Range("I46").Value = 694
$ws.Range("J46").Value = 4740
$ws.Range("K46").Value = 694
$ws.Range("L46").Value = 4740
$ws.Range("M46").Value = -506
$ws.Range("N46").Value = -5116

$ws.Range("H122").Value = 9133.733
$ws.Range("I122").Value = 35266.668
$ws.Range("J122").Value = 2600.5
$ws.Range("K122").Value = 105800.004
$ws.Range("L122").Value = 7801.5
$ws.Range("M122").Value = -103350.004
$ws.Range("N122").Value = -12701.5

$ws.Range("H132").Value = 43490176
$ws.Range("I132").Value = 71431720
$ws.Range("J132").Value = 25556.111
$ws.Range("K132").Value = 214295160
$ws.Range("L132").Value = 76668.333
$ws.Range("M132").Value = -214292630
$ws.Range("N132").Value = -81728.333

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 142864340
$ws.Range("I136").Value = 20416378
$ws.Range("J136").Value = 1000000000
$ws.Range("K136").Value = 61249134
$ws.Range("L136").Value = 3000000000
$ws.Range("M136").Value = -61246584
$ws.Range("N136").Value = -3000005100

$ws.Range("H137").Value = 50500
$ws.Range("J137").Value = 50500
$ws.Range("L137").Value = 50500
$ws.Range("N137").Value = -60700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26696.543
$ws.Range("I132").Value = 36014.2
$ws.Range("J132").Value = 9225.9375
$ws.Range("K132").Value = 108042.6
$ws.Range("L132").Value = 27677.8125
$ws.Range("M132").Value = -105512.6
$ws.Range("N132").Value = -32737.8125

$ws.Range("H136").Value = 4616
$ws.Range("I136").Value = 7257.0713
$ws.Range("J136").Value = 1658
$ws.Range("K136").Value = 21771.2139
$ws.Range("L136").Value = 4974
$ws.Range("M136").Value = -19221.2139
$ws.Range("N136").Value = -10074
